$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full set of rows 2-7 (header row 1 stays the same).
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics.
$rows = @(
  @{ Row = 2;  A = "FAPs"; B = "Bmp7"; C = "Acvr2a"; D = "ECs";
     E = 3; F = 1; G = 1.668521; H = 5.005563; I = 0.9677024783929865; J = 0.9677024783929865;
     K = 3; L = 1; M = 14.61878266666667; N = 43.856348; O = 0.2662829816142094; P = 0.2662829816142094;
     Q = 24.39174587376933; R = 219.525712863924; S = 0.2576827012619445; T = 0.2576827012619445 },

  @{ Row = 3;  A = "FAPs"; B = "Bmp7"; C = "Acvr2a"; D = "FAPs";
     E = 3; F = 1; G = 1.668521; H = 5.005563; I = 0.9677024783929865; J = 0.9677024783929865;
     K = 3; L = 1; M = 27.084169; N = 81.25250700000001; O = 0.4933415757187404; P = 0.4933415757187404;
     Q = 45.19050474404901; R = 406.7145426964411; S = 0.4774078655173263; T = 0.4774078655173263 },

  @{ Row = 4;  A = "FAPs"; B = "Bmp7"; C = "Acvr2a"; D = "sCs";
     E = 3; F = 1; G = 1.668521; H = 5.005563; I = 0.9677024783929865; J = 0.9677024783929865;
     K = 3; L = 1; M = 13.19647366666667; N = 39.589421; O = 0.2403754426670501; P = 0.2403754426670501;
     Q = 22.01859343878034; R = 198.167340949023; S = 0.2326119116137156; T = 0.2326119116137156 },

  @{ Row = 5;  A = "sCs"; B = "Bmp7"; C = "Acvr2a"; D = "ECs";
     E = 1; F = 0.3333333333333333; G = 0.05568766666666666; H = 0.167063; I = 0.03229752160701353; J = 0.03229752160701353;
     K = 3; L = 1; M = 14.61878266666667; N = 43.856348; O = 0.2662829816142094; P = 0.2662829816142094;
     Q = 0.8140858962137777; R = 7.326773065923999; S = 0.008600280352264917; T = 0.008600280352264917 },

  @{ Row = 6;  A = "sCs"; B = "Bmp7"; C = "Acvr2a"; D = "FAPs";
     E = 1; F = 0.3333333333333333; G = 0.05568766666666666; H = 0.167063; I = 0.03229752160701353; J = 0.03229752160701353;
     K = 3; L = 1; M = 27.084169; N = 81.25250700000001; O = 0.4933415757187404; P = 0.4933415757187404;
     Q = 1.508254175215667; R = 13.574287576941; S = 0.01593371020141412; T = 0.01593371020141412 },

  @{ Row = 7;  A = "sCs"; B = "Bmp7"; C = "Acvr2a"; D = "sCs";
     E = 1; F = 0.3333333333333333; G = 0.05568766666666666; H = 0.167063; I = 0.03229752160701353; J = 0.03229752160701353;
     K = 3; L = 1; M = 13.19647366666667; N = 39.589421; O = 0.2403754426670501; P = 0.2403754426670501;
     Q = 0.7348808267247778; R = 6.613927440523; S = 0.007763531053334495; T = 0.007763531053334494 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    $rowNum = $r["Row"]
    foreach ($col in $cols) {
        $cellAddr = "$col$rowNum"
        $ws.Range($cellAddr).Value = $r[$col]
    }
}

Write-Host "Updated rows 2-7 with new sCs cluster data"
